# Auto-generated Excel COM-interop script
# Applies numeric value updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# (market-price / profit recompute from a scheduled data-refresh run).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1372.5
$ws.Range("I28").Value = 1372.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1372.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -887.5
$ws.Range("N28").ClearContents()
$ws.Range("H62").Value = 3802.1538
$ws.Range("J62").Value = 3977.3333
$ws.Range("L62").Value = 3977.3333
$ws.Range("N62").Value = -5225.3333
$ws.Range("H65").Value = 3802.1538
$ws.Range("J65").Value = 3977.3333
$ws.Range("L65").Value = 19886.6665
$ws.Range("N65").Value = -26126.6665
$ws.Range("H69").Value = 29839.7
$ws.Range("H70").Value = 7776.2354
$ws.Range("J70").Value = 12098.6
$ws.Range("L70").Value = 36295.8
$ws.Range("N70").Value = -36835.8
$ws.Range("H72").Value = 29839.7
$ws.Range("H73").Value = 7776.2354
$ws.Range("J73").Value = 12098.6
$ws.Range("L73").Value = 36295.8
$ws.Range("N73").Value = -38167.8
$ws.Range("H74").Value = 9357.842000000001
$ws.Range("I74").Value = 7724.875
$ws.Range("J74").Value = 10545.454
$ws.Range("K74").Value = 7724.875
$ws.Range("L74").Value = 10545.454
$ws.Range("M74").Value = -6788.875
$ws.Range("N74").Value = -12417.454
$ws.Range("H77").Value = 9357.842000000001
$ws.Range("I77").Value = 7724.875
$ws.Range("J77").Value = 10545.454
$ws.Range("K77").Value = 38624.375
$ws.Range("L77").Value = 52727.27
$ws.Range("M77").Value = -33944.375
$ws.Range("N77").Value = -62087.27
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H112").Value = 1896.225
$ws.Range("J112").Value = 1922.3422
$ws.Range("L112").Value = 5767.0266
$ws.Range("N112").Value = -7983.0266
$ws.Range("H113").Value = 7621.8335
$ws.Range("I113").Value = 6024.467
$ws.Range("J113").Value = 9219.200000000001
$ws.Range("K113").Value = 6024.467
$ws.Range("L113").Value = 9219.200000000001
$ws.Range("M113").Value = -2770.467
$ws.Range("N113").Value = -15727.2
$ws.Range("H138").Value = 2819.3635
$ws.Range("I138").Value = 1452.68
$ws.Range("J138").Value = 3652.7073
$ws.Range("K138").Value = 4358.04
$ws.Range("L138").Value = 10958.1219
$ws.Range("M138").Value = 781.96
$ws.Range("N138").Value = -21238.1219

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23838.572
$ws.Range("J2").Value = 2887
$ws.Range("L2").Value = 2887
$ws.Range("N2").Value = -3113
$ws.Range("H37").Value = 995
$ws.Range("I37").Value = 995
$ws.Range("K37").Value = 995
$ws.Range("M37").Value = -722
$ws.Range("H61").Value = 3231.5557
$ws.Range("I61").Value = 3231.5557
$ws.Range("K61").Value = 3231.5557
$ws.Range("M61").Value = -3019.5557
$ws.Range("H92").Value = 62035
$ws.Range("J92").Value = 62035
$ws.Range("L92").Value = 62035
$ws.Range("N92").Value = -67027
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 23838.572
$ws.Range("J116").Value = 2887
$ws.Range("L116").Value = 2887
$ws.Range("N116").Value = -7475
$ws.Range("H136").Value = 3231.5557
$ws.Range("I136").Value = 3231.5557
$ws.Range("K136").Value = 9694.667099999999
$ws.Range("M136").Value = -7144.667099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23838.572
$ws.Range("J3").Value = 2887
$ws.Range("L3").Value = 2887
$ws.Range("N3").Value = -3115
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 65805.45
$ws.Range("J132").Value = 65805.45
$ws.Range("L132").Value = 65805.45
$ws.Range("N132").Value = -75925.45
$ws.Range("H134").Value = 2602.4167
$ws.Range("I134").Value = 1978.9524
$ws.Range("K134").Value = 5936.857199999999
$ws.Range("M134").Value = -3401.857199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2665.7144
$ws.Range("I58").Value = 1882
$ws.Range("K58").Value = 1882
$ws.Range("M58").Value = -1679
$ws.Range("H62").Value = 4312
$ws.Range("J62").Value = 3949
$ws.Range("L62").Value = 3949
$ws.Range("N62").Value = -5197
$ws.Range("H65").Value = 4312
$ws.Range("J65").Value = 3949
$ws.Range("L65").Value = 19745
$ws.Range("N65").Value = -25985
$ws.Range("H122").Value = 449030.88
$ws.Range("I122").Value = 732950.7
$ws.Range("J122").Value = 7377.778
$ws.Range("K122").Value = 2198852.1
$ws.Range("L122").Value = 22133.334
$ws.Range("M122").Value = -2196402.1
$ws.Range("N122").Value = -27033.334
$ws.Range("H134").Value = 7007.353
$ws.Range("I134").Value = 7131.9375
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 21395.8125
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -18860.8125
$ws.Range("N134").Value = -20112
$ws.Range("H136").Value = 2665.7144
$ws.Range("I136").Value = 1882
$ws.Range("K136").Value = 5646
$ws.Range("M136").Value = -3096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J2").Value = 99
$ws.Range("L2").Value = 594
$ws.Range("N2").Value = -820
$ws.Range("H107").Value = 2052.9443
$ws.Range("J107").Value = 2499.75
$ws.Range("L107").Value = 7499.25
$ws.Range("N107").Value = -11339.25
$ws.Range("H108").Value = 779
$ws.Range("I108").Value = 779
$ws.Range("K108").Value = 2337
$ws.Range("M108").Value = 543
$ws.Range("H122").Value = 4986.353
$ws.Range("I122").Value = 750.4
$ws.Range("J122").Value = 6751.3335
$ws.Range("K122").Value = 6753.599999999999
$ws.Range("L122").Value = 60762.0015
$ws.Range("M122").Value = -4303.599999999999
$ws.Range("N122").Value = -65662.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 99302
$ws.Range("J103").Value = 99302
$ws.Range("L103").Value = 99302
$ws.Range("N103").Value = -101646
$ws.Range("H113").Value = 9327.182000000001
$ws.Range("I113").Value = 3228.5715
$ws.Range("J113").Value = 19999.75
$ws.Range("K113").Value = 3228.5715
$ws.Range("L113").Value = 19999.75
$ws.Range("M113").Value = -1058.5715
$ws.Range("N113").Value = -24339.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3665.5
$ws.Range("I68").Value = 2886.8333
$ws.Range("K68").Value = 2886.8333
$ws.Range("M68").Value = -2137.8333
$ws.Range("H71").Value = 3665.5
$ws.Range("I71").Value = 2886.8333
$ws.Range("K71").Value = 14434.1665
$ws.Range("M71").Value = -10690.1665
$ws.Range("H82").Value = 1603.9048
$ws.Range("I82").Value = 1720.3334
$ws.Range("J82").Value = 1312.8334
$ws.Range("K82").Value = 1720.3334
$ws.Range("L82").Value = 1312.8334
$ws.Range("M82").Value = -1359.3334
$ws.Range("N82").Value = -2034.8334
$ws.Range("H85").Value = 1603.9048
$ws.Range("I85").Value = 1720.3334
$ws.Range("J85").Value = 1312.8334
$ws.Range("K85").Value = 1720.3334
$ws.Range("L85").Value = 1312.8334
$ws.Range("M85").Value = -472.3334
$ws.Range("N85").Value = -3808.8334
$ws.Range("H104").Value = 23933.334
$ws.Range("J104").Value = 23933.334
$ws.Range("L104").Value = 23933.334
$ws.Range("N104").Value = -30921.334
$ws.Range("H136").Value = 2795.2812
$ws.Range("I136").Value = 2565
$ws.Range("K136").Value = 7695
$ws.Range("M136").Value = -5145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 865.5
$ws.Range("I107").Value = 812.8
$ws.Range("J107").Value = 953.3333
$ws.Range("K107").Value = 2438.4
$ws.Range("L107").Value = 2859.9999
$ws.Range("M107").Value = -518.3999999999996
$ws.Range("N107").Value = -6699.9999
